$p = $ppt.ActivePresentation
$m = $p.SlideMaster
try {
  $m.Design.Name = "Office Theme"
  Write-Output "set via master.design.name ok -> $($m.Design.Name)"
} catch { Write-Output "err: $_" }

try {
  $m.Rename("Office Theme")
  Write-Output "rename ok"
} catch { Write-Output "err rename: $_" }
